$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title.
#    Built with InsertXML so the paragraph comes out with the same shape
#    (leading empty run, bold run, plain run) used elsewhere in this file.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter() | Out-Null

$metaPara = $d.Paragraphs.Item(2)
$metaRange = $d.Range($metaPara.Range.Start, $metaPara.Range.End)

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Discover what we like and don''t like about Book of Cats, a visually stunning slot from BGaming. Play for free and read our review.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaRange.InsertXML($metaXml) | Out-Null

# ---------------------------------------------------------------------------
# 2) Remove the duplicated bold "Play Book of Cats..." paragraph near the end
#    (it now only needs to live at the top of the document).
#    NOTE: Paragraph.Range.Text includes the trailing paragraph mark (chr 13),
#    so trim it before comparing against plain text.
# ---------------------------------------------------------------------------
$titleText = "Play Book of Cats for Free - Review of BGaming's Slot"
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if (($txt -eq $titleText) -and ($i -gt 2)) {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Replace the old meta-description text in the final (italic) paragraph
#    with the AI image-generation prompt.
# ---------------------------------------------------------------------------
$oldText = "Discover what we like and don't like about Book of Cats, a visually stunning slot from BGaming. Play for free and read our review."
$newText = 'Prompt: Create a cartoon-style feature image for the game "Book of Cats" that includes a happy Maya warrior with glasses. The Maya warrior should be holding a golden book in one hand and have a cat perched on their shoulder. The background of the image should feature palm trees and the Nile river with the sun shining brightly in the sky. The design should be colorful and eye-catching, with bold lines and a fun, playful style. The image should capture the magical atmosphere of the game and the theme of ancient Egyptian mythology.'

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq $oldText) {
        $target = $d.Range($p.Range.Start, $p.Range.End)
        $target.Text = $newText
        break
    }
}

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
